# Applies the "Parallelisierung der Feldprüfung" edit:
# Appends 16 new data rows (rows 5-20) to Tabelle1 (sheet1), repeating the
# existing part-number / article pattern with two new text values
# ("BARVERKAUF 1" and "test2"), and updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# The repeating 6-row data block (A, B, C) used for rows 5-10 and 11-16.
$block = @(
    @(10027, 50000, "BARVERKAUF 1"),
    @(10027, 70003, "test2"),
    @(10026, 50000, "BARVERKAUF 1"),
    @(10026, 70003, "test2"),
    @(10028, 50000, "BARVERKAUF"),
    @(10028, 70003, "test")
)

$row = 5

# First full repetition: rows 5-10
foreach ($entry in $block) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}

# Second full repetition: rows 11-16
foreach ($entry in $block) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}

# Final partial repetition (first four rows of the block): rows 17-20
for ($i = 0; $i -lt 4; $i++) {
    $entry = $block[$i]
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}

# Update the active selection on Tabelle1 to D14, matching the saved view.
$ws.Range("D14").Select()
